# Auto-generated edit script: update cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.968.28'
$ws.Range('E2').Value = '  +7.28%  '
$ws.Range('D3').Value = '2.453.34'
$ws.Range('E3').Value = '  +8.54%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '477.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +10.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +18.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.499'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +9.97%  '
$ws.Range('D9').Value = '2.449.01'
$ws.Range('E9').Value = '  +8.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0957'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +14.36%  '
$ws.Range('E11').Value = '  +6.01%  '
$ws.Range('E12').Value = '  +9.19%  '
$ws.Range('E13').Value = '  +1.98%  '
$ws.Range('D14').Value = '2.879.57'
$ws.Range('E14').Value = '  +8.82%  '
$ws.Range('D15').Value = '55.164.06'
$ws.Range('E15').Value = '  +7.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.36'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +11.48%  '
$ws.Range('E17').Value = '  +19.12%  '
$ws.Range('D18').Value = '2.453.96'
$ws.Range('E18').Value = '  +8.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.33'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +12.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.91'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +16.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '312.24'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.997'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('E23').Value = '  +13.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '57.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.64%  '
$ws.Range('E25').Value = '  +1.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.402'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.160'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +18.22%  '
$ws.Range('D28').Value = '2.555.94'
$ws.Range('E28').Value = '  +9.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.30'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +9.63%  '
$ws.Range('D30').Value = '0.0₃0766'
$ws.Range('E30').Value = '  +24.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '148.49'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.84'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +9.15%  '
$ws.Range('E34').Value = '  +13.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.13'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +12.46%  '
$ws.Range('E36').Value = '  +15.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.845'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +10.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.57'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.998'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '33.32'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.600'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.27%  '
$ws.Range('E42').Value = '  +12.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0539'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +10.58%  '
$ws.Range('E44').Value = '  +14.33%  '
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '254.57'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +33.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.61'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +17.94%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0890'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +11.74%  '
$ws.Range('D49').Value = '1.921.09'
$ws.Range('E49').Value = '  +3.19%  '
$ws.Range('E50').Value = '  +10.95%  '
$ws.Range('E51').Value = '  +11.69%  '
